$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force the cell to Text number format before assigning the value so
    # numeric-looking strings (e.g. "583.71", "1.00") are stored as text
    # instead of being auto-converted to numbers, matching the source data.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "63.177.39"
Set-TextValue "D3" "2.558.09"
Set-TextValue "E3" "  +0.75%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "583.71"
Set-TextValue "E5" "  +2.60%  "
Set-TextValue "D6" "147.99"
Set-TextValue "E6" "  +0.36%  "
Set-TextValue "E7" "  +0.05%  "
Set-TextValue "D8" "0.585"
Set-TextValue "E8" "  +0.93%  "
Set-TextValue "E9" "  +3.65%  "
Set-TextValue "D10" "5.61"
Set-TextValue "E10" "  -0.03%  "
Set-TextValue "E11" "  +0.26%  "
Set-TextValue "E12" "  +0.11%  "
Set-TextValue "D13" "27.61"
Set-TextValue "E13" "  +1.37%  "
Set-TextValue "D14" "3.014.46"
Set-TextValue "E14" "  +0.76%  "
Set-TextValue "D15" "62.998.55"
Set-TextValue "E15" "  +0.14%  "
Set-TextValue "E16" "  +4.01%  "
Set-TextValue "D17" "2.563.43"
Set-TextValue "E17" "  +0.62%  "
Set-TextValue "D18" "11.37"
Set-TextValue "E18" "  -0.92%  "
Set-TextValue "D19" "341.68"
Set-TextValue "E19" "  +2.36%  "
Set-TextValue "D20" "4.42"
Set-TextValue "E20" "  +3.22%  "
Set-TextValue "D21" "6.81"
Set-TextValue "E21" "  +1.41%  "
Set-TextValue "E22" "  -0.07%  "
Set-TextValue "D23" "66.34"
Set-TextValue "E23" "  +2.40%  "
Set-TextValue "E24" "  +3.35%  "
Set-TextValue "D25" "2.686.97"
Set-TextValue "E25" "  +0.56%  "
Set-TextValue "E26" "  +1.30%  "
Set-TextValue "D27" "8.11"
Set-TextValue "E27" "  +13.41%  "
Set-TextValue "E28" "  +2.00%  "
Set-TextValue "D31" "1.97"
Set-TextValue "E31" "  +5.72%  "
Set-TextValue "D32" "0.0₃0827"
Set-TextValue "E32" "  +1.50%  "
Set-TextValue "D33" "177.11"
Set-TextValue "E33" "  -0.15%  "
Set-TextValue "D34" "439.56"
Set-TextValue "E34" "  +6.32%  "
Set-TextValue "D35" "1.61"
Set-TextValue "E35" "  +2.49%  "
Set-TextValue "E36" "  +2.76%  "
Set-TextValue "D37" "19.26"
Set-TextValue "E37" "  +2.29%  "
Set-TextValue "D38" "4.51"
Set-TextValue "E38" "  +3.49%  "
Set-TextValue "E40" "  +0.20%  "
Set-TextValue "E41" "  +0.00%  "
Set-TextValue "D42" "150.84"
Set-TextValue "E42" "  -0.43%  "
Set-TextValue "E43" "  +2.39%  "
Set-TextValue "D44" "21.09"
Set-TextValue "E44" "  +2.44%  "
Set-TextValue "D45" "0.0550"
Set-TextValue "E45" "  +6.30%  "
Set-TextValue "D46" "0.611"
Set-TextValue "E46" "  +1.34%  "
Set-TextValue "E47" "  +1.31%  "
Set-TextValue "E48" "  +3.06%  "
Set-TextValue "D49" "18.40"
Set-TextValue "E49" "  +0.68%  "
Set-TextValue "E50" "  -1.95%  "
Set-TextValue "D51" "11.37"
Set-TextValue "E51" "  -0.28%  "

# Row 29/30 data swapped places (SuiNetwork now at 30, Binance-PegBSC-USD at 29)
Set-TextValue "B29" "Binance-PegBSC-USD"
Set-TextValue "C29" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.06%  "
Set-TextValue "B30" "SuiNetwork"
Set-TextValue "C30" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D30" "1.49"
Set-TextValue "E30" "  +0.33%  "
